# fix-reports: add header row for event-report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing placeholder row (was row 1) down to row 2 and
# insert a fresh blank row 1 above it for the human-readable header.
$ws.Rows.Item(1).Insert()

# Populate the new header row with Russian column captions.
$ws.Range("A1").Value = "Дата мероприятия"
$ws.Range("B1").Value = "Площадка"
$ws.Range("C1").Value = "Адресс площадки"
$ws.Range("D1").Value = "Кол-во билетов"
$ws.Range("E1").Value = "Сумма "
$ws.Range("F1").Value = "Кол-во заказов"

# The "EventReport" named range used by the report generator must now
# point at the placeholder row, which moved from row 1 to row 2.
$wb.Names.Item("EventReport").RefersTo = "=Лист1!`$A`$2:`$F`$2"

# Add a new named range covering the freshly added header row.
$wb.Names.Add("header", "=Лист1!`$A`$1:`$F`$1")

# Match the author's final selection (cell F2) before saving.
$ws.Range("F2").Select() | Out-Null
